# Updated cryptos list on Fri Dec 29 22:50:02 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '41.977.07'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -1.73%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.284.42'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -2.93%  '

# Row 4
$ws.Range('E4').Value = '  +0.06%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '311.53'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -3.94%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '105.99'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +3.17%  '

# Row 7
$ws.Range('E7').Value = '  -1.91%  '

# Row 8
$ws.Range('E8').Value = '  +0.08%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.607'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -2.80%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '40.42'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +0.69%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0909'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -1.27%  '

# Row 12
$ws.Range('E12').Value = '  -2.02%  '

# Row 13
$ws.Range('E13').Value = '  +0.01%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.967'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -3.12%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '15.46'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -4.37%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '2.628.43'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -2.98%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.307.63'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -1.83%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '41.923.66'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -1.60%  '

# Row 19
$ws.Range('E19').Value = '  -4.42%  '

# Row 20
$ws.Range('E20').Value = '  -1.94%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '73.33'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -4.19%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '3.44'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -6.43%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '255.94'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -2.85%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.32'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.35%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '9.32'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -7.14%  '

# Row 26
$ws.Range('E26').Value = '  +0.47%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '10.93'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -4.58%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.29'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +3.40%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '22.65'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.43%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '166.29'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -5.07%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '35.65'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +1.13%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0891'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -0.72%  '

# Row 33
$ws.Range('E33').Value = '  -5.86%  '

# Row 34
$ws.Range('E34').Value = '  -4.84%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '4.56'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +0.49%  '

# Row 38
$ws.Range('E38').Value = '  -1.27%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.79'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -0.40%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.63'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -4.05%  '

# Row 43
$ws.Range('E43').Value = '  -1.87%  '

# Row 44
$ws.Range('E44').Value = '  -4.31%  '

# Row 45
$ws.Range('E45').Value = '  +0.22%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '12.28'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +3.59%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '112.43'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -7.45%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '9.04'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -1.10%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '5.29'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -4.34%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '75.40'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +6.90%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.559.94'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +1.05%  '

# Row 35 (full content replace)
$ws.Range('B35').Value = 'Stellar'
$ws.Range('C35').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.129'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -2.16%  '

# Row 36 (full content replace)
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.118'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +7.50%  '

# Row 41 (full content replace)
$ws.Range('B41').Value = 'MultiversX'
$ws.Range('C41').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '71.37'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +2.19%  '

# Row 42 (full content replace)
$ws.Range('B42').Value = 'BitcoinSV'
$ws.Range('C42').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '97.99'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +5.58%  '

